$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = -2.6524319090109687
$ws.Range("F2").Value = 0.00006330237544809098
$ws.Range("G2").Value = 0.0015192570107541836
$ws.Range("H2").Value = 0.001465949747218949

$ws.Range("E3").Value = -2.2556893759404035
$ws.Range("F3").Value = 0.002462712123719609
$ws.Range("G3").Value = 0.029552545484635307
$ws.Range("H3").Value = 0.028515614064121786
